$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.199.94"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.063.81"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'248.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "'0.668"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D8").Value = "'57.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "'0.0787"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "'0.917"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.88%  "
$ws.Range("D14").Value = "2.361.06"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").Value = "2.061.62"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "'18.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +12.07%  "
$ws.Range("D18").Value = "37.211.84"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "'74.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").Value = "'238.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'2.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.63%  "
$ws.Range("D25").Value = "'9.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.67%  "
$ws.Range("D26").Value = "'2.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.20%  "
$ws.Range("D27").Value = "'170.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "'20.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").Value = "'5.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.15%  "
$ws.Range("D31").Value = "'1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("D32").Value = "'0.0627"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "'4.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.98%  "
$ws.Range("D34").Value = "'0.0892"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'2.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("E39").Value = "  +15.35%  "
$ws.Range("D40").Value = "'3.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.11%  "
$ws.Range("E41").Value = "  -10.73%  "
$ws.Range("D42").Value = "'17.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("D45").Value = "'96.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").Value = "1.277.11"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "'6.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").Value = "2.249.00"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "'44.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.18%  "
